$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.037.37"
$ws.Range("E2").Value = "  +0.86%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.747.08"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.85"
$ws.Range("E5").Value = "  +3.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9997"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5289"
$ws.Range("E7").Value = "  +2.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2781"
$ws.Range("E8").Value = "  +1.60%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06180"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.744.32"
$ws.Range("E10").Value = "  +0.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07228"
$ws.Range("E11").Value = "  +3.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.34"
$ws.Range("E12").Value = "  +1.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6438"
$ws.Range("E13").Value = "  +2.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.608"
$ws.Range("E14").Value = "  +2.82%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "78.26"
$ws.Range("E15").Value = "  +2.51%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9994"
$ws.Range("E16").Value = "  -0.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9991"
$ws.Range("E17").Value = "  -0.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.965.12"
$ws.Range("E18").Value = "  +0.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.59"
$ws.Range("E19").Value = "  +1.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000006730"
$ws.Range("E20").Value = "  +1.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.968.00"
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.316"
$ws.Range("E22").Value = "  +6.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.780"
$ws.Range("E23").Value = "  +4.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.203"
$ws.Range("E24").Value = "  +2.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "138.86"
$ws.Range("E25").Value = "  +1.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.516"
$ws.Range("E26").Value = "  +1.42%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.29"
$ws.Range("E27").Value = "  +2.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.805"
$ws.Range("E28").Value = "  -0.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "104.51"
$ws.Range("E29").Value = "  +1.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08303"
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.797"
$ws.Range("E31").Value = "  +5.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.683"
$ws.Range("E32").Value = "  +9.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04537"
$ws.Range("E33").Value = "  +3.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.644"
$ws.Range("E34").Value = "  +0.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9991"
$ws.Range("E35").Value = "  +3.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6325"
$ws.Range("E36").Value = "  +6.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.710"
$ws.Range("E37").Value = "  +1.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01593"
$ws.Range("E38").Value = "  +2.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.937"
$ws.Range("E39").Value = "  +0.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9991"
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.48"
$ws.Range("E41").Value = "  -2.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.3909"
$ws.Range("E42").Value = "  +2.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7408"
$ws.Range("E43").Value = "  +2.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.027"
$ws.Range("E44").Value = "  +3.57%  "
$ws.Range("E45").Value = "  +4.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.321"
$ws.Range("E46").Value = "  +2.91%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05346"
$ws.Range("E47").Value = "  -2.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.01"
$ws.Range("E48").Value = "  +4.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.70"
$ws.Range("E49").Value = "  +3.57%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.670"
$ws.Range("E50").Value = "  +3.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3454"
$ws.Range("E51").Value = "  +2.18%  "
